$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309; rows 309:351 shift down to 310:352.
$ws.Rows("309:309").Insert()

# Populate the newly inserted row 309 with the new weekly price record.
$ws.Range("A309").Value = 8
$ws.Range("B309").Value = "Terminal La Palmera de La Serena"
$ws.Range("C309").Value = "Coquimbo"
$ws.Range("D309").Value = 45142
$ws.Range("E309").Value = 4
$ws.Range("F309").Value = 100112037
$ws.Range("G309").Value = "Cebollín"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 1000
$ws.Range("K309").Value = 4500
$ws.Range("L309").Value = 5000
$ws.Range("M309").Value = 4750
$ws.Range("N309").Value = "$/paquete 36 unidades"
$ws.Range("O309").Value = "Provincia del Elquí"
$ws.Range("P309").Value = 132
$ws.Range("Q309").Value = 36
$ws.Range("R309").Value = "Hortaliza"
